$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells holding percentage-formatted text (e.g. "57%") need NumberFormat forced
# to Text ("@") before assignment, otherwise Excel auto-converts the literal
# "57%" string into the number 0.57 with a Percent number format, which changes
# the underlying cell type away from the original inline-string text value.

$ws.Range("E2").Value = "2026-02-09 18:18:31"
$ws.Range("E3").Value = "2026-02-09 18:18:34"
$ws.Range("E4").Value = "2026-02-09 18:18:36"
$ws.Range("E5").Value = "2026-02-09 18:18:39"
$ws.Range("O5").Value = "-2.8 °C"
$ws.Range("E6").Value = "2026-02-09 18:18:41"
$ws.Range("O6").Value = "8.8 °C"
$ws.Range("E7").Value = "2026-02-09 18:18:44"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "57%"
$ws.Range("E8").Value = "2026-02-09 18:18:47"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "68%"
$ws.Range("E9").Value = "2026-02-09 18:18:49"
$ws.Range("E10").Value = "2026-02-09 18:18:51"
$ws.Range("K10").Value = "11.0 MJ/m2"
$ws.Range("E11").Value = "2026-02-09 18:18:54"
$ws.Range("E12").Value = "2026-02-09 18:18:57"
$ws.Range("E13").Value = "2026-02-09 18:18:59"
$ws.Range("J13").Value = "1008.2 hPa"
$ws.Range("O13").Value = "2.9 °C"
$ws.Range("E14").Value = "2026-02-09 18:19:02"
$ws.Range("E15").Value = "2026-02-09 18:19:05"
$ws.Range("E16").Value = "2026-02-09 18:19:07"
$ws.Range("O16").Value = "-3.6 °C"
$ws.Range("E17").Value = "2026-02-09 18:19:10"
$ws.Range("O17").Value = "0.9 °C"
$ws.Range("E18").Value = "2026-02-09 18:19:13"
$ws.Range("E19").Value = "2026-02-09 18:19:15"
$ws.Range("E20").Value = "2026-02-09 18:19:18"
$ws.Range("O20").Value = "-4.5 °C"
$ws.Range("E21").Value = "2026-02-09 18:19:21"
$ws.Range("E22").Value = "2026-02-09 18:19:23"
$ws.Range("L22").Value = "33.1 km/h - 322º 17:35 TU"
$ws.Range("E23").Value = "2026-02-09 18:19:31"
$ws.Range("E24").Value = "2026-02-09 18:19:33"
$ws.Range("I24").Value = "0.3 mm"
$ws.Range("E25").Value = "2026-02-09 18:19:36"
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value = "73%"
$ws.Range("E26").Value = "2026-02-09 18:19:39"
$ws.Range("J26").Value = "1006.8 hPa"
$ws.Range("E27").Value = "2026-02-09 18:19:41"
$ws.Range("E28").Value = "2026-02-09 18:19:43"
$ws.Range("E29").Value = "2026-02-09 18:19:46"
$ws.Range("O29").Value = "8.8 °C"
$ws.Range("E30").Value = "2026-02-09 18:19:48"
$ws.Range("E31").Value = "2026-02-09 18:19:51"
$ws.Range("E32").Value = "2026-02-09 18:19:53"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "79%"
$ws.Range("O32").Value = "5.0 °C"
$ws.Range("E33").Value = "2026-02-09 18:19:56"
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = "77%"
$ws.Range("E34").Value = "2026-02-09 18:19:59"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = "73%"
$ws.Range("E35").Value = "2026-02-09 18:20:01"
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H35").Value = "76%"
$ws.Range("I35").Value = "0.3 mm"
$ws.Range("E36").Value = "2026-02-09 18:20:04"
$ws.Range("E37").Value = "2026-02-09 18:20:07"
$ws.Range("E38").Value = "2026-02-09 18:20:09"
$ws.Range("E39").Value = "2026-02-09 18:20:12"
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H39").Value = "77%"
$ws.Range("O39").Value = "-3.5 °C"
$ws.Range("E40").Value = "2026-02-09 18:20:14"
$ws.Range("E41").Value = "2026-02-09 18:20:17"
$ws.Range("E42").Value = "2026-02-09 18:20:19"
$ws.Range("E43").Value = "2026-02-09 18:20:22"
$ws.Range("E44").Value = "2026-02-09 18:20:25"
$ws.Range("E45").Value = "2026-02-09 18:20:27"
$ws.Range("J45").Value = "1007.4 hPa"
$ws.Range("E46").Value = "2026-02-09 18:20:30"
$ws.Range("O46").Value = "10.3 °C"
